$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update EC database: row 16 (Jorge Luis Covo Martinez) and row 17 (Paola Andrea Luengas
# Torres) swap places - each worker's document number, period and amounts (Valor Mora /
# Salario Basico) move together to the other row.
$ws.Range("C16").Value = "1102870612"
$ws.Range("D16").Value = "PAOLA ANDREA LUENGAS TORRES"
$ws.Range("E16").Value = "1907"
$ws.Range("F16").Value = 36000
$ws.Range("G16").Value = 955790

$ws.Range("C17").Value = "1129517708"
$ws.Range("D17").Value = "JORGE LUIS COVO MARTINEZ"
$ws.Range("E17").Value = "1908"
$ws.Range("F17").Value = 66250
$ws.Range("G17").Value = 1656232
